# The former row 22 ("ZLOZENIE" / 0 / "Centrum kompletacji" / "Kompletacja" / bordered
# blanks E:P) is removed entirely. This shifts the old row 23 ("P" / "Centrum
# kompletacji" / "Wydanie_na_produkcje") up to become the new row 22, and the old
# row 24 ("W" / "Wysylka" / "Wysylka") up to become the new row 23 - matching the
# target layout (sheet now ends at row 23 instead of row 24).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(22).Delete()

# Leave the sheet with the (now shifted-up) row selected, as Excel does after a
# row-header delete.
$ws.Rows.Item(22).Select() | Out-Null
